# Analyse des logs chgInfnegsRemet - Version 2
# Evolution de CompteursFichierColl.to_excel()
#
# 1. Rename the lone sheet "Sheet" -> "Feuille 1"
# 2. Add a new, empty sheet "Feuille 2" right after it (sheetId 2 / rId2)
# 3. Tag the header row (A1:C1) with the "highlight" named cell style
# 4. Refresh the timestamp stored in A2
# 5. Leave "Feuille 1" as the active / selected sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- rename the existing sheet --------------------------------------------
$ws.Name = "Feuille 1"

# -- header row gets the "highlight" cell style ----------------------------
$ws.Range("A1:C1").Style = "highlight"

# -- refresh the date/time value stored in A2 ------------------------------
$ws.Range("A2").Value = 43006.45625978009

# -- insert the new, empty worksheet right after "Feuille 1" --------------
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "Feuille 2"

# match the default page margins used throughout the workbook (in points:
# 0.75in/1in/0.5in -> 54/72/36)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# -- keep "Feuille 1" as the active tab ------------------------------------
$ws.Activate()
